$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header cell A1: "model" -> "Model"
$ws.Range("A1").Value = "Model"

# 2. Insert 5 new rows (alphabetically-sorted additions) and fill in their data.
#    Rows are inserted top-to-bottom at their FINAL target row numbers; because we
#    process them in ascending order, each insertion naturally lands in the right
#    spot relative to the previously-inserted rows.
#    NOTE: this interpreter only supports POSITIONAL function parameters.

function Insert-DataRow($RowNumber, $Name, $B, $C, $D, $E, $F) {
    # Insert a blank row, shifting this row and everything below it down by one.
    $ws.Rows.Item($RowNumber).Insert()

    # Copy formatting (borders/font/alignment/style) for column A from the row
    # directly above, so the new label cell matches the rest of the table (style "1").
    $ws.Range("A" + ($RowNumber - 1)).Copy()
    $ws.Range("A" + $RowNumber).PasteSpecial(-4122)

    $ws.Range("A" + $RowNumber).Value = $Name
    $ws.Range("B" + $RowNumber).Value = $B
    $ws.Range("C" + $RowNumber).Value = $C
    $ws.Range("D" + $RowNumber).Value = $D
    $ws.Range("E" + $RowNumber).Value = $E
    $ws.Range("F" + $RowNumber).Value = $F
}

Insert-DataRow 5 "centernet_resnet50_v2_512x512_coco17_tpu-8_bsize_16" 0.9909523809523809 0.990952380952381 0.988 0.9880000000000001 0.9693082883029307

Insert-DataRow 8 "efficientdet_d1_coco17_tpu-32" 0.0299010899010899 0.02931479402067638 0.04 0.0392156862745098 0.4228139897782832

Insert-DataRow 9 "efficientdet_d1_coco17_tpu-32_bsize_8" 0.9933333333333334 0.9933333333333334 0.992 0.9919999999999999 0.9641992746040781

Insert-DataRow 11 "faster_rcnn_resnet152_v1_640x640_coco17_tpu-8_bsize_8" 0.9626190476190477 0.9626190476190476 0.948 0.9479999999999998 0.9736932490821781

Insert-DataRow 21 "yolov8m_lego_416_bsize_8" 0.9933333333333334 0.9933333333333334 0.992 0.992 0.9610137022932269
